# agregadas nuevas tareas al scrum
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# --- Row 58: replace F58 ("pend") with E58 ("hecho") ---
$ws.Range("F58").ClearContents()
$ws.Range("E58").Value = "hecho"

# --- Row 74: fill in the previously empty B74 cell ---
$ws.Range("B74").Value = "investigar incrutacion de video"

# --- New rows 75-82: new scrum tasks, matching B74's formatting ---
$newTasks = @(
    "acomodar apk",
    "cortar y listar videos",
    "agregar videos al sistema",
    "hacer nuevas capturas para la interfaz",
    "completar marco teorico con info de discapacitados",
    "terminar conclusion y  lineas futuras",
    "hacer manual de usuario",
    "revision total final de la documentacion"
)

$row = 75
foreach ($task in $newTasks) {
    $cell = $ws.Range("B" + $row)
    $cell.Value = $task
    $cell.Font.Bold = $true
    $cell.VerticalAlignment = -4160
    $cell.WrapText = $true
    $row = $row + 1
}

# --- Update the active selection to match the new view state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 63
$ws.Range("B73").Select()
